$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the three new "Check defaults" rows (7, 8, 9)
$ws.Range("A7").Value = "Check defaults Page Profil"
$ws.Range("B7").Value = "<CHK>"
$ws.Range("C7").Value = "Check defaults"
$ws.Range("F7").Value = "<NOP>"

$ws.Range("A8").Value = "Check defaults Page Abwesenheiten"
$ws.Range("B8").Value = "<CHK>"
$ws.Range("D8").Value = "Check defaults"
$ws.Range("F8").Value = "<NOP>"

$ws.Range("A9").Value = "Check defaults Page Benachrichtigungen"
$ws.Range("B9").Value = "<CHK>"
$ws.Range("E9").Value = "Check defaults"
$ws.Range("F9").Value = "<NOP>"

# Reposition the picture (it was pushed down to make room for the new rows)
$shp = $ws.Shapes.Item(1)
$shp.Top = 152.4

# Update the selection on the sheet
$ws.Range("E7").Select()

$wb.Save()
